# Applies "tabelas atualizadas com novos dados" edits to the single
# table in the document:
#   - updates several statistic cells with revised values
#   - removes the CIRURGIA VASCULAR, PESQUISA CLINICA and REUMATOLOGIA rows

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Update individual cell values (row, column) -------------------------
# Row 2 = "n"
$t.Cell(2, 2).Range.Text = "325"
$t.Cell(2, 3).Range.Text = "86"

# Row 3 = "Grupo (%)" -> "test" column (NA -> 0.651)
$t.Cell(3, 4).Range.Text = "0.651"

# Row 7 = JOELHO
$t.Cell(7, 2).Range.Text = "113 (34.8)"
$t.Cell(7, 3).Range.Text = "36 (41.9)"

# Row 10 = OMBRO
$t.Cell(10, 3).Range.Text = "1 ( 1.2)"

# Row 11 = PÉ
$t.Cell(11, 3).Range.Text = "5 ( 5.8)"

# Row 13 = QUADRIL
$t.Cell(13, 2).Range.Text = "65 (20.0)"
$t.Cell(13, 3).Range.Text = "22 (25.6)"

# Row 15 = TRAUMA
$t.Cell(15, 2).Range.Text = "88 (27.1)"
$t.Cell(15, 3).Range.Text = "17 (19.8)"

# Row 16 = TUMOR
$t.Cell(16, 3).Range.Text = "3 ( 3.5)"

# --- Remove whole rows (delete from bottom to top to keep indices valid) -
# REUMATOLOGIA is row 14, PESQUISA CLINICA is row 12, CIRURGIA VASCULAR is row 4
$t.Rows.Item(14).Delete()
$t.Rows.Item(12).Delete()
$t.Rows.Item(4).Delete()
